# Add a new "2021" column (column R) to the 8.10.1 worksheet, mirroring
# the existing yearly columns (D..Q) that already live on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Values / formulas for the new column ---------------------------------
$ws.Range("R3").Value = 2021

$ws.Range("R4").Formula = "=R6/R8*100000"
$ws.Range("R5").Formula = "=R7/R8*100000"

$ws.Range("R6").Value = 312
$ws.Range("R7").Value = 1910
$ws.Range("R8").Value = 4409166

# --- Formatting: mirror column Q's formatting onto column R ---------------
$ws.Range("Q3:Q8").Copy()
$ws.Range("R3:R8").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- View: clear the old frozen/scrolled position and move the selection --
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("R15").Select()
